$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 38
$ws.Range("H38").Value = 3786.7
$ws.Range("J38").Value = 4724.875
$ws.Range("L38").Value = 14174.625
$ws.Range("N38").Value = -14918.625
# row 58
$ws.Range("H58").Value = 879.53845
$ws.Range("I58").Value = 241.7
$ws.Range("J58").Value = 3005.6667
$ws.Range("K58").Value = 725.0999999999999
$ws.Range("L58").Value = 9017.000100000001
$ws.Range("M58").Value = -575.0999999999999
$ws.Range("N58").Value = -9317.000100000001
# row 100
$ws.Range("H100").Value = 2278.5334
$ws.Range("I100").Value = 2272.4443
$ws.Range("K100").Value = 2272.4443
$ws.Range("M100").Value = -1731.4443
# row 132
$ws.Range("H132").Value = 5132384.5
$ws.Range("I132").Value = 7579597.5
$ws.Range("J132").Value = 4890.6665
$ws.Range("K132").Value = 22738792.5
$ws.Range("L132").Value = 14671.9995
$ws.Range("M132").Value = -22736262.5
$ws.Range("N132").Value = -19731.9995
# row 137
$ws.Range("H137").Value = 1173.7742
$ws.Range("I137").Value = 860.0857
$ws.Range("J137").Value = 1580.4073
$ws.Range("K137").Value = 2580.2571
$ws.Range("L137").Value = 4741.2219
$ws.Range("M137").Value = -30.25709999999981
$ws.Range("N137").Value = -9841.2219
# row 138
$ws.Range("H138").Value = 1202.49
$ws.Range("I138").Value = 661.8158
$ws.Range("J138").Value = 1533.871
$ws.Range("K138").Value = 1985.4474
$ws.Range("L138").Value = 4601.613
$ws.Range("M138").Value = 3154.5526
$ws.Range("N138").Value = -14881.613

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 4223.6196
$ws.Range("I32").Value = 3981.0176
$ws.Range("K32").Value = 3981.0176
$ws.Range("M32").Value = -3694.0176
# row 61
$ws.Range("H61").Value = 19231774
$ws.Range("I61").Value = 27027690
$ws.Range("J61").Value = 1843.4667
$ws.Range("K61").Value = 27027690
$ws.Range("L61").Value = 1843.4667
$ws.Range("M61").Value = -27027478
$ws.Range("N61").Value = -2267.4667
# row 74
$ws.Range("H74").Value = 1547.125
$ws.Range("I74").Value = 1144.6111
$ws.Range("J74").Value = 2754.6667
$ws.Range("K74").Value = 1144.6111
$ws.Range("L74").Value = 2754.6667
$ws.Range("M74").Value = -270.6111000000001
$ws.Range("N74").Value = -4502.6667
# row 77
$ws.Range("H77").Value = 1547.125
$ws.Range("I77").Value = 1144.6111
$ws.Range("J77").Value = 2754.6667
$ws.Range("K77").Value = 5723.0555
$ws.Range("L77").Value = 13773.3335
$ws.Range("M77").Value = -1355.0555
$ws.Range("N77").Value = -22509.3335
# row 81
$ws.Range("H81").Value = 49721.332
$ws.Range("I81").Value = 19164
$ws.Range("K81").Value = 19164
$ws.Range("M81").Value = -18166
# row 84
$ws.Range("H84").Value = 49721.332
$ws.Range("I84").Value = 19164
$ws.Range("K84").Value = 57492
$ws.Range("M84").Value = -52500
# row 132
$ws.Range("H132").Value = 1641.4
$ws.Range("I132").Value = 1432.76
$ws.Range("J132").Value = 1989.1333
$ws.Range("K132").Value = 4298.28
$ws.Range("L132").Value = 5967.3999
$ws.Range("M132").Value = -1768.28
$ws.Range("N132").Value = -11027.3999
# row 136
$ws.Range("H136").Value = 19231774
$ws.Range("I136").Value = 27027690
$ws.Range("J136").Value = 1843.4667
$ws.Range("K136").Value = 81083070
$ws.Range("L136").Value = 5530.4001
$ws.Range("M136").Value = -81080520
$ws.Range("N136").Value = -10630.4001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 1703.6129
$ws.Range("I20").Value = 1272.091
$ws.Range("J20").Value = 2758.4443
$ws.Range("K20").Value = 1272.091
$ws.Range("L20").Value = 2758.4443
$ws.Range("M20").Value = -1025.091
$ws.Range("N20").Value = -3252.4443
# row 105
$ws.Range("H105").Value = 112210870
$ws.Range("I105").Value = 144270930
$ws.Range("J105").Value = 654.5
$ws.Range("K105").Value = 144270930
$ws.Range("L105").Value = 654.5
$ws.Range("M105").Value = -144269183
$ws.Range("N105").Value = -4148.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2036.4324
$ws.Range("I31").Value = 1866.129
$ws.Range("K31").Value = 1866.129
$ws.Range("M31").Value = -1571.129
# row 34
$ws.Range("H34").Value = 2036.4324
$ws.Range("I34").Value = 1866.129
$ws.Range("K34").Value = 1866.129
$ws.Range("M34").Value = -1664.129
# row 86
$ws.Range("H86").Value = 1916122
$ws.Range("I86").Value = 2784248
$ws.Range("J86").Value = 22028.818
$ws.Range("K86").Value = 2784248
$ws.Range("L86").Value = 22028.818
$ws.Range("M86").Value = -2783125
$ws.Range("N86").Value = -24274.818
# row 89
$ws.Range("H89").Value = 1916122
$ws.Range("I89").Value = 2784248
$ws.Range("J89").Value = 22028.818
$ws.Range("K89").Value = 13921240
$ws.Range("L89").Value = 110144.09
$ws.Range("M89").Value = -13915624
$ws.Range("N89").Value = -121376.09
# row 99
$ws.Range("H99").Value = 2712.625
$ws.Range("I99").Value = 2636.6667
$ws.Range("J99").Value = 2940.5
$ws.Range("K99").Value = 2636.6667
$ws.Range("L99").Value = 2940.5
$ws.Range("M99").Value = -1138.6667
$ws.Range("N99").Value = -5936.5
# row 124
$ws.Range("H124").Value = 12740
$ws.Range("J124").Value = 12740
$ws.Range("L124").Value = 12740
$ws.Range("N124").Value = -17650
# row 126
$ws.Range("H126").Value = 2712.625
$ws.Range("I126").Value = 2636.6667
$ws.Range("J126").Value = 2940.5
$ws.Range("K126").Value = 7910.000100000001
$ws.Range("L126").Value = 8821.5
$ws.Range("M126").Value = -5440.000100000001
$ws.Range("N126").Value = -13761.5
# row 132
$ws.Range("H132").Value = 4251.7144
$ws.Range("I132").Value = 4923.613
$ws.Range("K132").Value = 14770.839
$ws.Range("M132").Value = -12240.839
# row 135
$ws.Range("H135").Value = 30275
$ws.Range("J135").Value = 33171.43
$ws.Range("L135").Value = 33171.43
$ws.Range("N135").Value = -43311.43

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 45
$ws.Range("H45").Value = 692.3333
$ws.Range("J45").Value = 692.3333
$ws.Range("L45").Value = 2076.9999
$ws.Range("N45").Value = -3140.9999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 86
$ws.Range("H86").Value = 28987.223
$ws.Range("J86").Value = 28987.223
$ws.Range("L86").Value = 28987.223
$ws.Range("N86").Value = -31359.223
# row 89
$ws.Range("H89").Value = 28987.223
$ws.Range("J89").Value = 28987.223
$ws.Range("L89").Value = 86961.66900000001
$ws.Range("N89").Value = -98817.66900000001
# row 102
$ws.Range("H102").Value = 20834240
$ws.Range("I102").Value = 35715124
$ws.Range("K102").Value = 35715124
$ws.Range("M102").Value = -35713502
# row 123
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
# row 126
$ws.Range("H126").Value = 2355.35
$ws.Range("I126").Value = 1686.2667
$ws.Range("K126").Value = 5058.800099999999
$ws.Range("M126").Value = -2588.800099999999
# row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
# row 132
$ws.Range("H132").Value = 2726.0344
$ws.Range("I132").Value = 2416.2273
$ws.Range("J132").Value = 3699.7144
$ws.Range("K132").Value = 7248.6819
$ws.Range("L132").Value = 11099.1432
$ws.Range("M132").Value = -4718.6819
$ws.Range("N132").Value = -16159.1432
# row 135
$ws.Range("H135").Value = 36666.668
$ws.Range("J135").Value = 34000
$ws.Range("L135").Value = 34000
$ws.Range("N135").Value = -44140

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 53
$ws.Range("H53").Value = 3250
$ws.Range("I53").Value = 3250
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 3250
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -2732
$ws.Range("N53").Value = $null
# row 136
$ws.Range("H136").Value = 4998.6
$ws.Range("I136").Value = 6143.591
$ws.Range("J136").Value = 1849.875
$ws.Range("K136").Value = 18430.773
$ws.Range("L136").Value = 5549.625
$ws.Range("M136").Value = -15880.773
$ws.Range("N136").Value = -10649.625

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 75
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872
# row 78
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360
# row 88
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
# row 91
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
# row 132
$ws.Range("H132").Value = 1549.4103
$ws.Range("I132").Value = 1524.4736
$ws.Range("J132").Value = 1573.1
$ws.Range("K132").Value = 4573.4208
$ws.Range("L132").Value = 4719.299999999999
$ws.Range("M132").Value = -2043.4208
$ws.Range("N132").Value = -9779.299999999999
# row 136
$ws.Range("H136").Value = 780.4
$ws.Range("I136").Value = 672.8182
$ws.Range("J136").Value = 1076.25
$ws.Range("K136").Value = 2018.4546
$ws.Range("L136").Value = 3228.75
$ws.Range("M136").Value = 531.5454
$ws.Range("N136").Value = -8328.75
